# Rename the "_old" / "_new" column-header suffixes to the respective
# format-version names ("_FV2410" / "_FV2504"), wrap the header/data range
# in an Excel Table (ListObject), and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1) cells -------------------------------
# Columns A:J carried the "_old" suffix -> "_FV2410"
$oldHeaders = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)
for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $oldHeaders[$i]
}

# Column K ("diff") is unchanged.

# Columns L:U carried the "_new" suffix -> "_FV2504"
$newHeaders = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, 12 + $i).Value = $newHeaders[$i]
}

# --- 2. Turn the used range into an Excel Table (ListObject) ----------
$lastRow = $ws.UsedRange.Rows.Count()
$lastCol = $ws.UsedRange.Columns.Count()
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$lo = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row ------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
